$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 18.01240722717635
$ws.Range("C2").Value = 7.053251811927651
$ws.Range("D2").Value = 13.93593755237917
$ws.Range("E2").Value = 14.4254752448334
$ws.Range("G2").Value = 3.725995161099738
$ws.Range("J2").Value = 8.644754712247629
$ws.Range("K2").Value = 13.63415955619782
$ws.Range("L2").Value = 12.33445245877133
$ws.Range("N2").Value = 22.67466171358504
$ws.Range("O2").Value = 33.76093678760819
# Row 3
$ws.Range("B3").Value = 17.87076515841797
$ws.Range("C3").Value = 7.019378704982575
$ws.Range("D3").Value = 13.93058913873882
$ws.Range("E3").Value = 14.44636800134981
$ws.Range("G3").Value = 3.728152806147501
$ws.Range("J3").Value = 8.65565721002149
$ws.Range("K3").Value = 13.53536402074438
$ws.Range("L3").Value = 12.34021383469987
$ws.Range("N3").Value = 22.73824987449331
$ws.Range("O3").Value = 33.83553556228343
# Row 4
$ws.Range("B4").Value = 17.78690706617667
$ws.Range("C4").Value = 6.998160192489275
$ws.Range("D4").Value = 13.92976497995518
$ws.Range("E4").Value = 14.4610931935967
$ws.Range("G4").Value = 3.729548719377811
$ws.Range("J4").Value = 8.662730066140876
$ws.Range("K4").Value = 13.47690989337503
$ws.Range("L4").Value = 12.34542360884598
$ws.Range("N4").Value = 22.77915322964724
$ws.Range("O4").Value = 33.88694672501256
# Row 5
$ws.Range("B5").Value = 17.75354904494469
$ws.Range("C5").Value = 6.98940975588437
$ws.Range("D5").Value = 13.93004932501577
$ws.Range("E5").Value = 14.46757121918707
$ws.Range("G5").Value = 3.730135503496933
$ws.Range("J5").Value = 8.665707817029952
$ws.Range("K5").Value = 13.45366582594336
$ws.Range("L5").Value = 12.34796796354521
$ws.Range("N5").Value = 22.79629067359728
$ws.Range("O5").Value = 33.90930594464203
# Row 6
$ws.Range("B6").Value = 17.74806010986088
$ws.Range("C6").Value = 6.987950545442466
$ws.Range("D6").Value = 13.93013404133145
$ws.Range("E6").Value = 14.4686757339263
$ws.Range("G6").Value = 3.730234023645115
$ws.Range("J6").Value = 8.666208046865581
$ws.Range("K6").Value = 13.44984159466277
$ws.Range("L6").Value = 12.34841592111669
$ws.Range("N6").Value = 22.79916469723199
$ws.Range("O6").Value = 33.91310372297145
# Row 7
$ws.Range("B7").Value = 17.78645384698943
$ws.Range("C7").Value = 6.998042598405134
$ws.Range("D7").Value = 13.9297663016007
$ws.Range("E7").Value = 14.46117862510205
$ws.Range("G7").Value = 3.729556560247185
$ws.Range("J7").Value = 8.662769838059472
$ws.Range("K7").Value = 13.4765940538997
$ws.Range("L7").Value = 12.34545621593864
$ws.Range("N7").Value = 22.77938245056623
$ws.Range("O7").Value = 33.88724256701979
# Row 8
$ws.Range("B8").Value = 17.96294340649895
$ws.Range("C8").Value = 7.04166076969181
$ws.Range("D8").Value = 13.93358397333215
$ws.Range("E8").Value = 14.43228559839286
$ws.Range("G8").Value = 3.72672439240315
$ws.Range("J8").Value = 8.64843546914001
$ws.Range("K8").Value = 13.59965006712828
$ws.Range("L8").Value = 12.33609242400196
$ws.Range("N8").Value = 22.69620165564572
$ws.Range("O8").Value = 33.78549396611058
# Row 9
$ws.Range("B9").Value = 18.33222460197624
$ws.Range("C9").Value = 7.123776711817865
$ws.Range("D9").Value = 13.96050346289983
$ws.Range("E9").Value = 14.39066004535291
$ws.Range("G9").Value = 3.721732157993463
$ws.Range("J9").Value = 8.623317613274658
$ws.Range("K9").Value = 13.85748268328446
$ws.Range("L9").Value = 12.33095971275656
$ws.Range("N9").Value = 22.54778245270364
$ws.Range("O9").Value = 33.63050765185349
# Row 10
$ws.Range("B10").Value = 18.61558910559912
$ws.Range("C10").Value = 7.181912709182988
$ws.Range("D10").Value = 13.99199856400087
$ws.Range("E10").Value = 14.36921704494666
$ws.Range("G10").Value = 3.718403129163793
$ws.Range("J10").Value = 8.606669511018053
$ws.Range("K10").Value = 14.05560650516482
$ws.Range("L10").Value = 12.33519727503919
$ws.Range("N10").Value = 22.44761388150152
$ws.Range("O10").Value = 33.54385921149626
# Row 11
$ws.Range("B11").Value = 18.74665623717404
$ws.Range("C11").Value = 7.20786050414185
$ws.Range("D11").Value = 14.00883685015584
$ws.Range("E11").Value = 14.36144011528707
$ws.Range("G11").Value = 3.716961457060874
$ws.Range("J11").Value = 8.599484214562118
$ws.Range("K11").Value = 14.14732072404773
$ws.Range("L11").Value = 12.33884987478558
$ws.Range("N11").Value = 22.40395355169048
$ws.Range("O11").Value = 33.51036252621113
# Row 12
$ws.Range("B12").Value = 18.79655961335138
$ws.Range("C12").Value = 7.2176125429357
$ws.Range("D12").Value = 14.01557072014123
$ws.Range("E12").Value = 14.35877891631928
$ws.Range("G12").Value = 3.716425931911419
$ws.Range("J12").Value = 8.596818836400059
$ws.Range("K12").Value = 14.18225241708441
$ws.Range("L12").Value = 12.34047964648433
$ws.Range("N12").Value = 22.38769339118121
$ws.Range("O12").Value = 33.49853016463218
# Row 13
$ws.Range("B13").Value = 18.7858006097562
$ws.Range("C13").Value = 7.21551558429209
$ws.Range("D13").Value = 14.01410461657943
$ws.Range("E13").Value = 14.3593394435782
$ws.Range("G13").Value = 3.716540804906087
$ws.Range("J13").Value = 8.597390407301598
$ws.Range("K13").Value = 14.17472071641073
$ws.Range("L13").Value = 12.34011770184694
$ws.Range("N13").Value = 22.39118318140083
$ws.Range("O13").Value = 33.50104056685232
# Row 14
$ws.Range("B14").Value = 18.75075661442279
$ws.Range("C14").Value = 7.208664302055354
$ws.Range("D14").Value = 14.00938370451245
$ws.Range("E14").Value = 14.36121549373957
$ws.Range("G14").Value = 3.716917190855065
$ws.Range("J14").Value = 8.599263820706117
$ws.Range("K14").Value = 14.15019069034049
$ws.Range("L14").Value = 12.33897902321347
$ws.Range("N14").Value = 22.40261035250512
$ws.Range("O14").Value = 33.50937198925909
# Row 15
$ws.Range("B15").Value = 18.72932527049814
$ws.Range("C15").Value = 7.20445801609604
$ws.Range("D15").Value = 14.00653846721877
$ws.Range("E15").Value = 14.36240156249426
$ws.Range("G15").Value = 3.717149091583777
$ws.Range("J15").Value = 8.600418565548765
$ws.Range("K15").Value = 14.13519078872118
$ws.Range("L15").Value = 12.33831362140628
$ws.Range("N15").Value = 22.4096453513252
$ws.Range("O15").Value = 33.51458621621703
# Row 16
$ws.Range("B16").Value = 18.60706367850532
$ws.Range("C16").Value = 7.180206744621448
$ws.Range("D16").Value = 13.99094837802582
$ws.Range("E16").Value = 14.36976502828209
$ws.Range("G16").Value = 3.718498804659767
$ws.Range("J16").Value = 8.607146872702874
$ws.Range("K16").Value = 14.04964242031806
$ws.Range("L16").Value = 12.3349931591688
$ws.Range("N16").Value = 22.45050545939896
$ws.Range("O16").Value = 33.5461674989402
# Row 17
$ws.Range("B17").Value = 18.53258584834449
$ws.Range("C17").Value = 7.165200158611476
$ws.Range("D17").Value = 13.98202510029438
$ws.Range("E17").Value = 14.37478836446634
$ws.Range("G17").Value = 3.719345397904942
$ws.Range("J17").Value = 8.611373665617217
$ws.Range("K17").Value = 13.99754877124613
$ws.Range("L17").Value = 12.33339697416303
$ws.Range("N17").Value = 22.47605938519152
$ws.Range("O17").Value = 33.56705850182778
# Row 18
$ws.Range("B18").Value = 18.4899541951361
$ws.Range("C18").Value = 7.156521967374019
$ws.Range("D18").Value = 13.9771292446334
$ws.Range("E18").Value = 14.37786383656927
$ws.Range("G18").Value = 3.719839183776026
$ws.Range("J18").Value = 8.613841341582603
$ws.Range("K18").Value = 13.96773684397185
$ws.Range("L18").Value = 12.33264134032667
$ws.Range("N18").Value = 22.49093687727356
$ws.Range("O18").Value = 33.57963164706081
# Row 19
$ws.Range("B19").Value = 18.47555644837358
$ws.Range("C19").Value = 7.153575708250262
$ws.Range("D19").Value = 13.97551232602888
$ws.Range("E19").Value = 14.3789371334156
$ws.Range("G19").Value = 3.720007548907593
$ws.Range("J19").Value = 8.614683137322501
$ws.Range("K19").Value = 13.95766977377057
$ws.Range("L19").Value = 12.33241343821803
$ws.Range("N19").Value = 22.49600501342987
$ws.Range("O19").Value = 33.58398436930957
# Row 20
$ws.Range("B20").Value = 18.54049310645836
$ws.Range("C20").Value = 7.166802494430772
$ws.Range("D20").Value = 13.98295053792023
$ws.Range("E20").Value = 14.37423435692004
$ws.Range("G20").Value = 3.719254568263438
$ws.Range("J20").Value = 8.610919936780347
$ws.Range("K20").Value = 14.00307879750287
$ws.Range("L20").Value = 12.33355008784929
$ws.Range("N20").Value = 22.4733205517484
$ws.Range("O20").Value = 33.56477694516947
# Row 21
$ws.Range("B21").Value = 18.76104284735513
$ws.Range("C21").Value = 7.210678710356877
$ws.Range("D21").Value = 14.01076067578183
$ws.Range("E21").Value = 14.36065675654385
$ws.Range("G21").Value = 3.71680635520938
$ws.Range("J21").Value = 8.598712048541952
$ws.Range("K21").Value = 14.15739049818554
$ws.Range("L21").Value = 12.33930679972024
$ws.Range("N21").Value = 22.39924651367285
$ws.Range("O21").Value = 33.50690171710509
# Row 22
$ws.Range("B22").Value = 18.90674637758352
$ws.Range("C22").Value = 7.238923046045163
$ws.Range("D22").Value = 14.03101859156619
$ws.Range("E22").Value = 14.35343662198566
$ws.Range("G22").Value = 3.715266929576754
$ws.Range("J22").Value = 8.591057085961399
$ws.Range("K22").Value = 14.25940417956563
$ws.Range("L22").Value = 12.34450584037628
$ws.Range("N22").Value = 22.35242592154905
$ws.Range("O22").Value = 33.47404388516961
# Row 23
$ws.Range("B23").Value = 18.82885202927726
$ws.Range("C23").Value = 7.223888679498389
$ws.Range("D23").Value = 14.0200172609142
$ws.Range("E23").Value = 14.35713905925036
$ws.Range("G23").Value = 3.716083020137567
$ws.Range("J23").Value = 8.595113158814369
$ws.Range("K23").Value = 14.20486008253956
$ws.Range("L23").Value = 12.34160005826627
$ws.Range("N23").Value = 22.3772697442579
$ws.Range("O23").Value = 33.49112603048074
# Row 24
$ws.Range("B24").Value = 18.53691765063657
$ws.Range("C24").Value = 7.166078236195071
$ws.Range("D24").Value = 13.98253141763736
$ws.Range("E24").Value = 14.3744842396811
$ws.Range("G24").Value = 3.719295610324313
$ws.Range("J24").Value = 8.611124950334306
$ws.Range("K24").Value = 14.00057824258992
$ws.Range("L24").Value = 12.33348036032757
$ws.Range("N24").Value = 22.47455819823343
$ws.Range("O24").Value = 33.56580668440417
# Row 25
$ws.Range("B25").Value = 18.23005323573271
$ws.Range("C25").Value = 7.101939595907461
$ws.Range("D25").Value = 13.95115270530832
$ws.Range("E25").Value = 14.40031383189543
$ws.Range("G25").Value = 3.723022938604518
$ws.Range("J25").Value = 8.629794224247805
$ws.Range("K25").Value = 13.78610310193655
$ws.Range("L25").Value = 12.33093732930181
$ws.Range("N25").Value = 22.58636911573084
$ws.Range("O25").Value = 33.66765961337871
